$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 0.0003937
$ws.Range("F2").Value = 0.01766242
$ws.Range("G2").Value = 0.0007100866538461538

$ws.Range("E3").Value = 0.650059925
$ws.Range("F3").Value = 0.67256134
$ws.Range("G3").Value = 0.659487631
